$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row appended to the rate history table.
# Column A holds a date-shaped string ("2025-09-06"); force text formatting
# first so Excel doesn't auto-convert it into a date serial number, keeping
# it a literal string like the other cells in the sheet.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2025-09-06"
$ws.Range("B7").Value = "15:16:31"
$ws.Range("C7").Value = "1.00 EUR = 1614.4992 ARS"
